# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型" sheets
# to reflect the newer scrape snapshot, per commit:
#   "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Row -> new value for the "展览" sheet (column F)
$sheet1Updates = @{
    3  = 5324
    5  = 53
    7  = 602
    8  = 570
    9  = 1050
    11 = 1464
    12 = 4263
    13 = 438
    14 = 187
    15 = 165
    17 = 3387
    18 = 166
    19 = 1084
    20 = 102
    22 = 197
    24 = 39
    27 = 306
    28 = 30
    31 = 24
    32 = 26
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# Row -> new value for the "全部类型" sheet (column F)
$sheet4Updates = @{
    4  = 5324
    6  = 53
    8  = 602
    9  = 570
    10 = 1050
    12 = 1464
    13 = 4263
    14 = 438
    15 = 187
    16 = 165
    18 = 3387
    19 = 166
    20 = 1084
    21 = 102
    23 = 197
    25 = 39
    28 = 306
    29 = 30
    32 = 24
    33 = 26
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
